$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new article row (row 2): name, article number, GTIN
$ws.Range("A2").Value = "AMD Ryzen 9 9900X3D"
$ws.Range("B2").Value = 106861
$ws.Range("C2").Value = 730143315579

# Move the active selection to C5 (matches the saved selection in the sheet view)
$ws.Range("C5").Select()
